$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old "edit discussion" note from B27
$ws.Range("B27").ClearContents()

# Move the "working on it" note into B30 (replacing former "edit discussion" content/slot)
# Set this first so it takes shared-string index 33, matching the target string order.
$ws.Range("B30").Value = "working on it with new sim_learnLH_pea_optimal.m"

# Add a note about objective_CB_approx.m in B28
$ws.Range("B28").Value = "Note: objective_CB_approx.m should use the learning code " + [char]0x2026 + "univariate.m, but there's a problem there, so I don't do that yet. 30 July 2020"

# Add a new note in B31
$ws.Range("B31").Value = "need a number there too"

# Update the selected cell to B32, matching the post-edit cursor location
$ws.Range("B32").Select()
